$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15, shifting existing rows 15-54 down to 16-55.
$ws.Rows.Item(15).Insert()

# Fill in the new row 15 with the new data record.
# Columns that are constant across every data row in this sheet:
$ws.Range("A15").Value = 2
$ws.Range("B15").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C15").Value = "Coquimbo"
$ws.Range("D15").Value = 44742
$ws.Range("E15").Value = 4
$ws.Range("F15").Value = 100112026
$ws.Range("G15").Value = "Haba"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 600
$ws.Range("K15").Value = 13000
$ws.Range("L15").Value = 15000
$ws.Range("M15").Value = 14000
$ws.Range("N15").Value = "$/saco 25 kilos"
$ws.Range("O15").Value = "Provincia de Limarí"
$ws.Range("P15").Value = 560
$ws.Range("Q15").Value = 25
$ws.Range("R15").Value = "Hortaliza"

# Ensure the date cell keeps the date number format used by the rest of column D.
$ws.Range("D15").NumberFormat = $ws.Range("D16").NumberFormat
